$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Income")

# New income row: "Giant Marketers" salary (base + bonus entered as an
# equation string), posted to Dutch Bangla Bank.
$ws.Range("A3").Value = "Giant Marketers"
$ws.Range("B3").Value = "Salary"

# C3/D3 look numeric/date-like ("25000+5000", "2026-02-01"); a leading
# apostrophe keeps them stored as literal text instead of being
# auto-converted to a number/date by Excel's input parsing.
$ws.Range("C3").Value = "'25000+5000"
$ws.Range("D3").Value = "'2026-02-01"

$ws.Range("E3").Value = "Dutch Bangla Bank"
$ws.Range("F3").Value = "me"

# G3/I3/J3 stay blank (notes / recurringFrequency / recurringNextDate are
# unused for this entry) but the cells themselves still need to exist, so
# format them as text before clearing them to an empty string.
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = ""

$ws.Range("H3").Value = $false

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = ""

$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = ""

# K3 (paymentMethod) is intentionally left untouched/absent, matching the
# source row which has no paymentMethod recorded.
